# meetingnote2023.1.2.docx - "Updated meeting note with more comments"
#
# 1) Add a new bulleted "Comments" item (same numbered list, numId 18) right
#    after the last screenshot / "The minimum size..." bullet and before the
#    "Specific tasks before next meeting." heading, followed by a blank
#    spacer paragraph.
# 2) Mark a page break as falling right before "Less urgent tasks".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 2 first: it only rewrites an existing paragraph's content (no
# paragraph-count shift), so doing it before the insertions below keeps
# every other paragraph index stable regardless of ordering.
# ---------------------------------------------------------------------
$lessUrgent = $d.Paragraphs.Item(27)
if ($lessUrgent.Range.Text.Trim() -ne "Less urgent tasks") {
    throw "Unexpected paragraph 27 content: $($lessUrgent.Range.Text)"
}
$lessUrgentXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space='preserve'>      </w:t></w:r>" +
    "<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t xml:space='preserve'>Less urgent tasks </w:t></w:r>" +
    "</w:p>"
$lessUrgent.Range.InsertXML($lessUrgentXml)

# ---------------------------------------------------------------------
# Change 1: insert the two new paragraphs after the screenshot paragraph
# that follows "The minimum size should not be there anymore."
# ---------------------------------------------------------------------
$shotPara = $d.Paragraphs.Item(18)
if ($shotPara.Range.InlineShapes.Count -ne 1) {
    throw "Unexpected paragraph 18 content - expected the screenshot paragraph"
}

# New list-item paragraph with the feedback comment.
$shotPara.Range.InsertParagraphAfter()
$commentPara = $d.Paragraphs.Item(19)
$commentXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='18'/></w:numPr>" +
    "<w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>User</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>&#8217;s</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> (Dr. Jiang</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>&#8217;s test during the meeting</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>) survey feedback was not saved to database.</w:t></w:r>" +
    "</w:p>"
$commentPara.Range.InsertXML($commentXml)

# Blank spacer paragraph right after it (matches the style used elsewhere
# between sections: left-indent 720 / hanging 360, 12pt run formatting).
$commentPara = $d.Paragraphs.Item(19)
$commentPara.Range.InsertParagraphAfter()
$spacerPara = $d.Paragraphs.Item(20)
$spacerXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:pPr><w:ind w:left='720' w:hanging='360'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "</w:p>"
$spacerPara.Range.InsertXML($spacerXml)

Write-Output "Done"
